# Applies the weekly FlashScore odds refresh:
#  - updates several odds values in row 2 (Tigre - Instituto)
#  - updates one odds value in row 5 (Botafogo RJ - Vitoria)
#  - removes the old "ECUADOR - LIGA PRO" fixture (row 10), which shifts
#    the "USA - MLS" fixture (old row 11) up into row 10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Tigre - Instituto) odds updates ---
$ws.Range("G2").Value = 2.4
$ws.Range("I2").Value = 3.4
$ws.Range("J2").Value = 3.2
$ws.Range("L2").Value = 4
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5
$ws.Range("X2").Value = 10
$ws.Range("Z2").Value = 23
$ws.Range("AA2").Value = 23
$ws.Range("AK2").Value = 34
$ws.Range("AN2").Value = 4.33
$ws.Range("AO2").Value = 15
$ws.Range("BB2").Value = 301

# --- Row 5 (Botafogo RJ - Vitoria) odds update ---
$ws.Range("N5").Value = 13

# --- Remove the old row 10 (ECUADOR - LIGA PRO fixture) ---
# This shifts the following row (USA - MLS, old row 11) up to become the
# new row 10, matching the expected final sheet shape (A1:BD10).
$ws.Rows.Item(10).Delete()
